$d = $word.ActiveDocument

# 1. Shorten the "Debounced persistence" bullet in section 3 (LOCKED DESIGN DECISIONS)
$d.Content.Find.Execute(
    "Debounced persistence with flush on close (and flush before any action that can close the dialog, e.g., sheet activation)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Debounced persistence with flush on close", 2) | Out-Null

# 2. Remove the UI-layout-invariants / Favorites-Recents-allocator / Recents-always-shown bullets
#    (section 3), keeping the trailing line break that already followed them.
$d.Content.Find.Execute(
    "• UI layout invariants must be centralized (e.g., listbox containment + row-label truncation) via shared primitives/helpers; avoid duplicated inline styles^l• Favorites/Recents layout uses a two-scenario allocator:^l  - No-conflict: show all Favorites and Recents (subject to minimum 20% share each); any extra space is placed between sections^l  - Conflict: apply user-selected policy (fixed ratio with surplus-donation, or prioritize Favorites up to 80%)^l• Recents section is always shown (even when empty) for feature discoverability and to preserve consistent tab structure^l",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "", 2) | Out-Null

# 3. Remove the "Decision Recording Rule" bullet and its sub-bullets (section 4)
$d.Content.Find.Execute(
    "• Decision Recording Rule (Required): Any non-trivial UX/architecture/workflow decision must be recorded in this LPD with:^l  - Decision^l  - Rationale^l  - Implications / tradeoffs^l  - Revisit conditions^l",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "", 2) | Out-Null

# 4. Remove the "Dialog containers can scroll unexpectedly..." bullet (section 6)
$d.Content.Find.Execute(
    "^l• Dialog containers can scroll unexpectedly unless the root/container overflow is explicitly locked; enforce overflow: hidden at the dialog level and allow scrolling only inside listboxes",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "", 2) | Out-Null

# 5. Remove the "Make dialog height..." bullet (section 7)
$d.Content.Find.Execute(
    "^l• Make dialog height choose a sensible size at open time based on available viewport; listboxes should scale accordingly while preserving " + [char]0x201C + "no dialog scrolling" + [char]0x201D + " (scroll only inside listboxes)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "", 2) | Out-Null

# 6. Append the new "Decision Log Rules" / "Decisions Added (2026-01-25)" sections,
#    preceded by a page break, at the end of the document body.
$end = $d.Content
$end.Collapse(0)
$end.InsertXML(@'
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:br w:type='page'/></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='Heading1'/></w:pPr><w:r><w:t>Decision Log Rules</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>Rule: Every non-trivial decision (UX, architecture, persistence, performance, workflow) must be recorded with its rationale.</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>Format per decision: Decision • Rationale • Implications • Revisit Conditions.</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='Heading1'/></w:pPr><w:r><w:t>Decisions Added (2026-01-25)</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>Decision: Store Recents as a bounded history (MAX_RECENTS) without excluding the active sheet; apply exclusions only when rendering the Recents list.</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>Rationale: Separates persistence (history) from UI affordances (jump targets). Avoids the confusing N−1 behavior where a user setting of N displayed only N−1 due to post-slice filtering.</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>Implications: UI rendering filters out active sheet and non-displayable sheets first, then slices to the user-configured count. Storage remains stable and future-proof for analytics/history views.</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>Revisit Conditions: If Recents should reflect manual Excel navigation or require multi-window state, revisit where history is collected and how active exclusion is defined.</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>Decision: Use a single shared constant MAX_RECENTS as the source of truth for both storage cap and Settings UI maximum.</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>Rationale: Prevents silent truncation and drift where UI allows N but storage only retains &lt;N. Reduces maintenance risk by avoiding duplicated constants.</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>Implications: Shared constants module imported by both dialog and storage layers; any change to MAX_RECENTS updates both automatically.</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>Revisit Conditions: If we intentionally want different limits for stored vs displayed recents, split into MAX_RECENTS_STORED and MAX_RECENTS_DISPLAY (still in the same shared constants module).</w:t></w:r></w:p>
'@)
